# "implementada a interface da tv e corrigidos alguns bugs"
#
# Sheet1 held 3 rows of "appliance" data (name / label / value / on-off).
# Row 1 ("TV da sala") loses its extra boolean column and its custom row
# height; it is replaced by a new first row for the living-room A/C unit.
# The old "lamp sala" row is replaced by the TV row (now carrying a
# wattage + an extra numeric column), and the old "Ar quarto" row is
# removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D carried left-over custom widths/styles from the previous
# layout; reset that formatting before (re)writing the cells so nothing
# bleeds into the new data (e.g. column C's numeric style onto the new C1).
$ws.Columns("A:D").ClearFormats()

# Row 1: "ar da sala" / A/C / 23 / FALSE
#   (was: "TV da sala" / Televisor / <blank> / TRUE / FALSE)
$ws.Range("A1").Value = "ar da sala"
$ws.Range("B1").Value = "A/C"
$ws.Range("C1").Value = 23
$ws.Range("D1").Value = $false
$ws.Range("E1").ClearContents()

# Row 2: "TV da sala" / Televisor / 12 / 36 / TRUE
#   (was: "lamp sala" / Lampada / 0 / FALSE)
$ws.Range("A2").Value = "TV da sala"
$ws.Range("B2").Value = "Televisor"
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 36
$ws.Range("E2").Value = $true

# Row 3 ("Ar quarto" / A/C / 19 / FALSE) no longer exists.
$ws.Range("A3:E3").ClearContents()

# Row 1 no longer has the old custom row height (18.75pt); let it go back
# to the sheet default.
$ws.Rows.Item(1).AutoFit()
